$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 63 was blank; fill in the missing work-log entry (date, activity, hours)
# so that the C82 total picks it up correctly (bug: row wasn't counted).
$ws.Range("A63").Value = 43227
$ws.Range("B63").Value = "Commentaires"
$ws.Range("C63").Value = 4

# Update the current selection/scroll position left behind by the editor.
[void]$ws.Range("E57").Select()
